$d = $word.ActiveDocument

# Set the paragraph-mark language first so it folds into w:pPr/w:rPr/w:lang.
$p = $d.Paragraphs(1)
$p.Range.LanguageID = "en-US"

# Insert the new text at the very start of the document (before the bookmark).
$r = $d.Range(0, 0)
$r.Text = "rhtxsjfuyjtjdufyjdytfjhty"

# Mark the newly inserted run with the same language.
$d.Range(0, 25).LanguageID = "en-US"
